$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Extend the "Dirigirse al profesor" acceptance clause with the new
#    schedule text, right before the final period.
# ---------------------------------------------------------------------------
$r1 = $d.Content
$r1.Find.Execute(
    "si no escribe con ese formato el mensaje será ignorado.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "si no escribe con ese formato el mensaje será ignorado en horario de 9:00 AM a 5:00 PM de lunes a viernes.",
    2
)

# Move the hidden "_GoBack" bookmark: delete it from its old spot (after
# "Dirigirse al profesor" near the end of the document) and recreate it
# right after the period that closes the sentence we just edited.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$r2 = $d.Content
$r2.Find.Execute(
    "de lunes a viernes.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "", 0
)
$period = $d.Range($r2.End - 1, $r2.End)
$d.Bookmarks.Add("_GoBack", $period)

# ---------------------------------------------------------------------------
# 2. Merge "El examen no será corregido " + "y perderá los puntos" into a
#    single run.
# ---------------------------------------------------------------------------
$r3 = $d.Content
$r3.Find.Execute(
    "El examen no será corregido y perderá los puntos",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "El examen no será corregido y perderá los puntos",
    2
)

# ---------------------------------------------------------------------------
# 3. Merge "Semana 1" + "4" + ": " into a single run "Semana 14: ".
# ---------------------------------------------------------------------------
$r4 = $d.Content
$r4.Find.Execute(
    "Semana 14: ",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Semana 14: ",
    2
)
